# chore: update Sheets via scheduled runner
# Refreshes cached market-price / profit figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) on a handful of leve rows across several
# job sheets, matching the latest market snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1884.3334
$ws.Range("I40").Value = 1708.4286
$ws.Range("K40").Value = 1708.4286
$ws.Range("M40").Value = -1533.4286
$ws.Range("H55").Value = 540.5
$ws.Range("I55").Value = 130
$ws.Range("K55").Value = 130
$ws.Range("M55").Value = 84
$ws.Range("H106").Value = 22996.666
$ws.Range("I106").Value = 27662.223
$ws.Range("K106").Value = 27662.223
$ws.Range("M106").Value = -27031.223
$ws.Range("H132").Value = 3203.8
$ws.Range("I132").Value = 3029.75
$ws.Range("J132").Value = 3900
$ws.Range("K132").Value = 9089.25
$ws.Range("L132").Value = 11700
$ws.Range("M132").Value = -6559.25
$ws.Range("N132").Value = -16760
$ws.Range("H137").Value = 3015.0454
$ws.Range("I137").Value = 1276.2727
$ws.Range("J137").Value = 4753.8184
$ws.Range("K137").Value = 3828.8181
$ws.Range("L137").Value = 14261.4552
$ws.Range("M137").Value = -1278.8181
$ws.Range("N137").Value = -19361.4552
$ws.Range("H138").Value = 4825.1665
$ws.Range("I138").Value = 4666
$ws.Range("K138").Value = 13998
$ws.Range("M138").Value = -8858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2012.5
$ws.Range("I45").Value = 2014.2858
$ws.Range("K45").Value = 2014.2858
$ws.Range("M45").Value = -1637.2858
$ws.Range("H61").Value = 2065.6667
$ws.Range("I61").Value = 2065.6667
$ws.Range("K61").Value = 2065.6667
$ws.Range("M61").Value = -1853.6667
$ws.Range("H74").Value = 2109.2778
$ws.Range("I74").Value = 1062.2727
$ws.Range("K74").Value = 1062.2727
$ws.Range("M74").Value = -188.2727
$ws.Range("H77").Value = 2109.2778
$ws.Range("I77").Value = 1062.2727
$ws.Range("K77").Value = 5311.363499999999
$ws.Range("M77").Value = -943.3634999999995
$ws.Range("H102").Value = 1949
$ws.Range("I102").Value = 1949
$ws.Range("K102").Value = 1949
$ws.Range("M102").Value = -327
$ws.Range("H122").Value = 4269.9375
$ws.Range("I122").Value = 3881.9
$ws.Range("K122").Value = 11645.7
$ws.Range("M122").Value = -9195.700000000001
$ws.Range("H132").Value = 1562.4717
$ws.Range("I132").Value = 1492.2449
$ws.Range("K132").Value = 4476.7347
$ws.Range("M132").Value = -1946.7347
$ws.Range("H136").Value = 2065.6667
$ws.Range("I136").Value = 2065.6667
$ws.Range("K136").Value = 6197.000100000001
$ws.Range("M136").Value = -3647.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 26198.8
$ws.Range("J88").Value = 26198.8
$ws.Range("L88").Value = 26198.8
$ws.Range("N88").Value = -27010.8
$ws.Range("H91").Value = 26198.8
$ws.Range("J91").Value = 26198.8
$ws.Range("L91").Value = 26198.8
$ws.Range("N91").Value = -29006.8
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H134").Value = 974.34283
$ws.Range("I134").Value = 709.5
$ws.Range("J134").Value = 3799.3333
$ws.Range("K134").Value = 2128.5
$ws.Range("L134").Value = 11397.9999
$ws.Range("M134").Value = 406.5
$ws.Range("N134").Value = -16467.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1720.85
$ws.Range("I134").Value = 1225.4333
$ws.Range("J134").Value = 3207.1
$ws.Range("K134").Value = 3676.2999
$ws.Range("L134").Value = 9621.299999999999
$ws.Range("M134").Value = -1141.2999
$ws.Range("N134").Value = -14691.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1805.8462
$ws.Range("I34").Value = 1362.2858
$ws.Range("K34").Value = 4086.8574
$ws.Range("M34").Value = -4002.8574
$ws.Range("H68").Value = 1525.1818
$ws.Range("I68").Value = 1291.6
$ws.Range("J68").Value = 1719.8334
$ws.Range("K68").Value = 3874.8
$ws.Range("L68").Value = 5159.5002
$ws.Range("M68").Value = -3063.8
$ws.Range("N68").Value = -6781.5002
$ws.Range("H71").Value = 1525.1818
$ws.Range("I71").Value = 1291.6
$ws.Range("J71").Value = 1719.8334
$ws.Range("K71").Value = 11624.4
$ws.Range("L71").Value = 15478.5006
$ws.Range("M71").Value = -7568.4
$ws.Range("N71").Value = -23590.5006
$ws.Range("H129").Value = 3515.4546
$ws.Range("I129").Value = 4899.25
$ws.Range("J129").Value = 2724.7144
$ws.Range("K129").Value = 14697.75
$ws.Range("L129").Value = 8174.1432
$ws.Range("M129").Value = -9697.75
$ws.Range("N129").Value = -18174.1432
$ws.Range("H134").Value = 3407.5715
$ws.Range("J134").Value = 14999
$ws.Range("L134").Value = 44997
$ws.Range("N134").Value = -55137
$ws.Range("H139").Value = 3396.087
$ws.Range("I139").Value = 2722.4285
$ws.Range("J139").Value = 4444
$ws.Range("K139").Value = 8167.2855
$ws.Range("L139").Value = 13332
$ws.Range("M139").Value = -3027.2855
$ws.Range("N139").Value = -23612

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 281.53845
$ws.Range("I2").Value = 76.75
$ws.Range("K2").Value = 76.75
$ws.Range("M2").Value = 36.25
$ws.Range("H102").Value = 1221.6578
$ws.Range("I102").Value = 647.2
$ws.Range("J102").Value = 2326.3845
$ws.Range("K102").Value = 647.2
$ws.Range("L102").Value = 2326.3845
$ws.Range("M102").Value = 974.8
$ws.Range("N102").Value = -5570.3845
$ws.Range("H105").Value = 20814.334
$ws.Range("J105").Value = 20814.334
$ws.Range("L105").Value = 20814.334
$ws.Range("N105").Value = -27802.334
$ws.Range("H122").Value = 581352.2
$ws.Range("I122").Value = 72979
$ws.Range("K122").Value = 218937
$ws.Range("M122").Value = -216487
$ws.Range("H132").Value = 2372.8845
$ws.Range("I132").Value = 1883
$ws.Range("K132").Value = 5649
$ws.Range("M132").Value = -3119

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2809.0588
$ws.Range("I132").Value = 1890.7307
$ws.Range("K132").Value = 5672.1921
$ws.Range("M132").Value = -3142.1921
$ws.Range("H136").Value = 2871.5
$ws.Range("I136").Value = 2871.5
$ws.Range("K136").Value = 8614.5
$ws.Range("M136").Value = -6064.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1070.15
$ws.Range("I81").Value = 1105.421
$ws.Range("J81").Value = 400
$ws.Range("K81").Value = 2210.842
$ws.Range("L81").Value = 800
$ws.Range("M81").Value = -1149.842
$ws.Range("N81").Value = -2922
$ws.Range("H84").Value = 1070.15
$ws.Range("I84").Value = 1105.421
$ws.Range("J84").Value = 400
$ws.Range("K84").Value = 11054.21
$ws.Range("L84").Value = 4000
$ws.Range("M84").Value = -5750.210000000001
$ws.Range("N84").Value = -14608
$ws.Range("H140").Value = 95000
$ws.Range("J140").Value = 95000
$ws.Range("L140").Value = 95000
$ws.Range("N140").Value = -105360
